# Refactor the water-types table:
#  - drop the "pH" column (column B)
#  - drop the "Sample" header label, turning that column into a
#    zero-based row index (values shift from 1..10 down to 0..9)
#  - insert a "-" separator into the generated "Code" column values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B ("pH"); this shifts MCation/MAnion/Salinity/Alkalinity/
# BEX/Code one column to the left (C:H -> B:G) while column A ("Sample")
# keeps its values, style and position for now.
$ws.Columns.Item(2).Delete()

# The old "Sample" header text is dropped entirely (no header label for
# the new index column), so fully clear A1 (value + style).
$ws.Range("A1").Clear()

# Re-number the index column: old Sample values were 1..10, new values
# are 0-based (0..9).
for ($r = 2; $r -le 11; $r++) {
  $ws.Cells.Item($r, 1).Value = $r - 2
}

# Update the generated Code column (now column G) so it contains a "-"
# between the salinity/alkalinity prefix and the facies name.
$codes = @("g*-MgSO4", "g*-NaCl", "g*-CaMIX", "g*-CaHCO3", "g3-CaHCO3+", "F2-CaMIX+", "F2-CaNO3+", "F2-CaHCO3", "F3-CaMIX+", "F1-CaNO3+")
for ($r = 2; $r -le 11; $r++) {
  $ws.Cells.Item($r, 7).Value = $codes[$r - 2]
}
